# Append: 2026-01-08 06:39 JST
# Update the "取得日時" (acquisition datetime) column A on the first
# worksheet ("ランサーズ") for every existing data row (rows 2-6),
# replacing the old timestamp with the new run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-08 06:39:45"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
